$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2038.7778
$ws.Range("J17").Value = 2038.7778
$ws.Range("L17").Value = 6116.3334
$ws.Range("N17").Value = -6452.3334
$ws.Range("H18").Value = 2483.3333
$ws.Range("I18").Value = 1225
$ws.Range("J18").Value = 5000
$ws.Range("K18").Value = 1225
$ws.Range("L18").Value = 5000
$ws.Range("M18").Value = -941
$ws.Range("N18").Value = -5568
$ws.Range("H112").Value = 3632.8333
$ws.Range("I112").Value = 2549.5
$ws.Range("K112").Value = 7648.5
$ws.Range("M112").Value = -6540.5
$ws.Range("H132").Value = 292422.84
$ws.Range("J132").Value = 16335.167
$ws.Range("L132").Value = 49005.501
$ws.Range("N132").Value = -54065.501
$ws.Range("H134").Value = 115196
$ws.Range("J134").Value = 115196
$ws.Range("L134").Value = 115196
$ws.Range("N134").Value = -125336
$ws.Range("H136").Value = 79202.414
$ws.Range("J136").Value = 100061.29
$ws.Range("L136").Value = 100061.29
$ws.Range("N136").Value = -110261.29
$ws.Range("H138").Value = 1902.11
$ws.Range("I138").Value = 1321.875
$ws.Range("J138").Value = 1952.5652
$ws.Range("K138").Value = 3965.625
$ws.Range("L138").Value = 5857.6956
$ws.Range("M138").Value = 1174.375
$ws.Range("N138").Value = -16137.6956
$ws.Range("H141").Value = 2839.6956
$ws.Range("I141").Value = 2032.2632
$ws.Range("J141").Value = 6675
$ws.Range("K141").Value = 6096.7896
$ws.Range("L141").Value = 20025
$ws.Range("M141").Value = -916.7896000000001
$ws.Range("N141").Value = -30385

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 775.6531
$ws.Range("I2").Value = 657.8919
$ws.Range("K2").Value = 657.8919
$ws.Range("M2").Value = -544.8919
$ws.Range("H32").Value = 2537.13
$ws.Range("I32").Value = 2537.13
$ws.Range("K32").Value = 2537.13
$ws.Range("M32").Value = -2250.13
$ws.Range("H45").Value = 3833.6667
$ws.Range("I45").Value = 3449.5833
$ws.Range("J45").Value = 4345.778
$ws.Range("K45").Value = 3449.5833
$ws.Range("L45").Value = 4345.778
$ws.Range("M45").Value = -3072.5833
$ws.Range("N45").Value = -5099.778
$ws.Range("H61").Value = 3341.7334
$ws.Range("I61").Value = 2939.1292
$ws.Range("K61").Value = 2939.1292
$ws.Range("M61").Value = -2727.1292
$ws.Range("H112").Value = 84999.5
$ws.Range("J112").Value = 84999.5
$ws.Range("L112").Value = 84999.5
$ws.Range("N112").Value = -87953.5
$ws.Range("H116").Value = 775.6531
$ws.Range("I116").Value = 657.8919
$ws.Range("K116").Value = 657.8919
$ws.Range("M116").Value = 1636.1081
$ws.Range("H132").Value = 19439.074
$ws.Range("I132").Value = 23471.334
$ws.Range("J132").Value = 11064.385
$ws.Range("K132").Value = 70414.00199999999
$ws.Range("L132").Value = 33193.155
$ws.Range("M132").Value = -67884.00199999999
$ws.Range("N132").Value = -38253.155
$ws.Range("H136").Value = 3341.7334
$ws.Range("I136").Value = 2939.1292
$ws.Range("K136").Value = 8817.3876
$ws.Range("M136").Value = -6267.3876

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 775.6531
$ws.Range("I3").Value = 657.8919
$ws.Range("K3").Value = 657.8919
$ws.Range("M3").Value = -543.8919
$ws.Range("H122").Value = 78499
$ws.Range("J122").Value = 78499
$ws.Range("L122").Value = 78499
$ws.Range("N122").Value = -88299
$ws.Range("H123").Value = 20000
$ws.Range("I123").Value = 20000
$ws.Range("K123").Value = 20000
$ws.Range("M123").Value = -15100
$ws.Range("H134").Value = 2545.4
$ws.Range("I134").Value = 2476.4634
$ws.Range("K134").Value = 7429.3902
$ws.Range("M134").Value = -4894.3902

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1067.421
$ws.Range("I107").Value = 965.7692
$ws.Range("J107").Value = 1287.6666
$ws.Range("K107").Value = 965.7692
$ws.Range("L107").Value = 1287.6666
$ws.Range("M107").Value = 954.2308
$ws.Range("N107").Value = -5127.6666
$ws.Range("H132").Value = 31009808
$ws.Range("I132").Value = 35089524
$ws.Range("K132").Value = 105268572
$ws.Range("M132").Value = -105266042
$ws.Range("H134").Value = 2079
$ws.Range("I134").Value = 2021
$ws.Range("K134").Value = 6063
$ws.Range("M134").Value = -3528

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 759.9375
$ws.Range("I5").Value = 474.92307
$ws.Range("J5").Value = 1995
$ws.Range("K5").Value = 1424.76921
$ws.Range("L5").Value = 5985
$ws.Range("M5").Value = -1312.76921
$ws.Range("N5").Value = -6209
$ws.Range("H7").Value = 361.35715
$ws.Range("I7").Value = 208.42857
$ws.Range("K7").Value = 625.28571
$ws.Range("M7").Value = -513.28571
$ws.Range("H135").Value = 759.9375
$ws.Range("I135").Value = 474.92307
$ws.Range("J135").Value = 1995
$ws.Range("K135").Value = 4274.30763
$ws.Range("L135").Value = 17955
$ws.Range("M135").Value = -1739.30763
$ws.Range("N135").Value = -23025
$ws.Range("H136").Value = 479132.2
$ws.Range("I136").Value = 589987.5
$ws.Range("J136").Value = 7997
$ws.Range("K136").Value = 1769962.5
$ws.Range("L136").Value = 23991
$ws.Range("M136").Value = -1764862.5
$ws.Range("N136").Value = -34191
$ws.Range("H138").Value = 95654.73
$ws.Range("I138").Value = 116015.11
$ws.Range("K138").Value = 348045.33
$ws.Range("M138").Value = -342905.33
$ws.Range("H140").Value = 2771.8845
$ws.Range("I140").Value = 2090
$ws.Range("K140").Value = 6270
$ws.Range("M140").Value = -1090

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1534.7646
$ws.Range("I113").Value = 1512.0714
$ws.Range("J113").Value = 1640.6666
$ws.Range("K113").Value = 1512.0714
$ws.Range("L113").Value = 1640.6666
$ws.Range("M113").Value = 657.9286
$ws.Range("N113").Value = -5980.6666
$ws.Range("H132").Value = 96927.62
$ws.Range("I132").Value = 125592.875
$ws.Range("J132").Value = 5198.8
$ws.Range("K132").Value = 376778.625
$ws.Range("L132").Value = 15596.4
$ws.Range("M132").Value = -374248.625
$ws.Range("N132").Value = -20656.4

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 20005496
$ws.Range("I40").Value = 33337032
$ws.Range("K40").Value = 33337032
$ws.Range("M40").Value = -33336896
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("H132").Value = 5589.9707
$ws.Range("I132").Value = 4906.38
$ws.Range("J132").Value = 7488.8335
$ws.Range("K132").Value = 14719.14
$ws.Range("L132").Value = 22466.5005
$ws.Range("M132").Value = -12189.14
$ws.Range("N132").Value = -27526.5005
$ws.Range("H133").Value = 87200
$ws.Range("J133").Value = 87200
$ws.Range("L133").Value = 87200
$ws.Range("N133").Value = -92260
$ws.Range("H136").Value = 3105.963
$ws.Range("I136").Value = 2289.6086
$ws.Range("K136").Value = 6868.825800000001
$ws.Range("M136").Value = -4318.825800000001
$ws.Range("M100").ClearContents()

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7248688
$ws.Range("I132").Value = 12347336
$ws.Range("J132").Value = 3242.0527
$ws.Range("K132").Value = 37042008
$ws.Range("L132").Value = 9726.158100000001
$ws.Range("M132").Value = -37039478
$ws.Range("N132").Value = -14786.1581
